$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.877306342124939
$ws.Range("B1").Value = 3.058123588562012
$ws.Range("C1").Value = 2.618117094039917
$ws.Range("D1").Value = 2.265565395355225
$ws.Range("E1").Value = 1.663646340370178
